$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Main sheet ("Main" / sheet1): new "SoundHound" label, updated Price,
# and the "CEO: Keyvan Mohajer" note moving from row 9 to row 10.
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Main")

$wsMain.Range("B2").Value = "SoundHound"
$wsMain.Range("M2").Value = 4.61

$wsMain.Range("L9").ClearContents()
$wsMain.Range("L10").Value = "CEO: Keyvan Mohajer"

# ---------------------------------------------------------------------
# Model sheet (sheet2): extend the year header row out to 2034 and
# convert the existing per-quarter formula rows into fill-right shared
# formulas (matches what Excel produces when you fill a formula across
# a row).
# ---------------------------------------------------------------------
$wsModel = $wb.Worksheets.Item("Model")

$wsModel.Range("R2").Value = 2021
$wsModel.Range("S2").Formula = "=+R2+1"
$wsModel.Range("T2:AE2").Formula = "=+S2+1"

$wsModel.Range("E5:L5").Formula = "=+E3-E4"
$wsModel.Range("E9:L9").Formula = "=+E8+E7+E6"
$wsModel.Range("E10:L10").Formula = "=+E5-E9"

# ---------------------------------------------------------------------
# View state: zoom levels + active selections on each sheet, leaving
# "Main" as the selected/active tab at the end (matches tabSelected="1"
# staying on sheet1 in the saved file).
# ---------------------------------------------------------------------
$wsModel.Select()
$excel.ActiveWindow.Zoom = 145

$wsMain.Select()
$excel.ActiveWindow.Zoom = 190
$wsMain.Range("A3").Select()
